$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row from the sheet's used range
$lastRow = $ws.UsedRange.Rows.Count

# Swap the contents of columns C and D (codeforiati:group-code <-> codeforiati:group-name)
# for every row, including the header row.
for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
